$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of "lesson 2" data.
# Shared-string insertion order matters: "uree-unupytu.wav" must become a
# shared string before "unupytu.wav" (B16 is written first, then B15).
$ws.Range("B16").Value = "uree-unupytu.wav"
$ws.Range("A16").Value = 2

$ws.Range("B15").Value = "unupytu.wav"
$ws.Range("A15").Value = 2

# Selection ends on D16 after the edits.
$ws.Range("D16").Select()

# Page setup: orientation explicitly set to portrait (adds <pageSetup .../>).
$ws.PageSetup.Orientation = 1
